# Auto-generated edit script: Add data for 2024-02-18
# Updates the 2024 (column K) and a few 2023 (column J) values across the
# "Citywide Totals", "By Neighborhood" and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 884
$ws.Range("K3").Value = 834
$ws.Range("J4").Value = 1792
$ws.Range("K4").Value = 192
$ws.Range("K5").Value = 49
$ws.Range("I6").Value = 8963
$ws.Range("K6").Value = 1155
$ws.Range("I7").Value = 26236
$ws.Range("J7").Value = 29251
$ws.Range("K7").Value = 3114

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 49
$ws.Range("K3").Value = 57
$ws.Range("I6").Value = 495
$ws.Range("K6").Value = 64
$ws.Range("I7").Value = 1541
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 47
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 93
$ws.Range("I8").Value = 1541
$ws.Range("K8").Value = 187
$ws.Range("K11").Value = 69
$ws.Range("K15").Value = 20
$ws.Range("K19").Value = 83
$ws.Range("J20").Value = 635
$ws.Range("K20").Value = 71
$ws.Range("K23").Value = 29
$ws.Range("K27").Value = 38
$ws.Range("K29").Value = 153
$ws.Range("K30").Value = 6
$ws.Range("K31").Value = 34
$ws.Range("K33").Value = 132
$ws.Range("K36").Value = 37
$ws.Range("K37").Value = 94
$ws.Range("K42").Value = 100
$ws.Range("K47").Value = 21
$ws.Range("K51").Value = 46
$ws.Range("J52").Value = 742
$ws.Range("K54").Value = 53
$ws.Range("J55").Value = 460
$ws.Range("K56").Value = 4
$ws.Range("K60").Value = 23
$ws.Range("J63").Value = 85
$ws.Range("K65").Value = 87
$ws.Range("K67").Value = 135
$ws.Range("K69").Value = 9
$ws.Range("K71").Value = 10
$ws.Range("K72").Value = 13
$ws.Range("K77").Value = 20
$ws.Range("K78").Value = 41
$ws.Range("K79").Value = 85
$ws.Range("K84").Value = 27
$ws.Range("K85").Value = 154
$ws.Range("K90").Value = 28
$ws.Range("K92").Value = 13
$ws.Range("K93").Value = 10
$ws.Range("K96").Value = 48
$ws.Range("K97").Value = 24
$ws.Range("K98").Value = 19
$ws.Range("I101").Value = 26236
$ws.Range("J101").Value = 29251
$ws.Range("K101").Value = 3114

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K3").Value = 4
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 37
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 45
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 26
$ws.Range("K3").Value = 27
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 21
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 41

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 460

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 27
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 15
$ws.Range("J4").Value = 53
$ws.Range("J7").Value = 635
$ws.Range("K7").Value = 71

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 9
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 14
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 11
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K4").Value = 2
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 12
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 48
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K4").Value = 2
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 4

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 209
$ws.Range("J7").Value = 742
